# Actualiza base de datos EC y agrega parte 1 de nuevos estado de cuenta
#
# The "Periodo Mora" rows (16 and 17) need their period/value data updated:
#  - Row 16 now corresponds to period 2402 (valor mora 38133)
#  - Row 17 now corresponds to period 2403 (valor mora 52000)
#  - "Salario Basico" (column G) is updated from 1300000 to 1423500 for both rows

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Hoja1")

# Row 16: period 2402
$ws.Range("E16").Value = "2402"
$ws.Range("F16").Value = 38133
$ws.Range("G16").Value = 1423500

# Row 17: period 2403
$ws.Range("E17").Value = "2403"
$ws.Range("F17").Value = 52000
$ws.Range("G17").Value = 1423500
